$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.25
$ws.Range("C2").Value = 0.2
$ws.Range("D2").Value = 0.2222222222222222
$ws.Range("C3").Value = 0.5714285714285714
$ws.Range("D3").Value = 0.5333333333333333
$ws.Range("B5").Value = 0.375
$ws.Range("C5").Value = 0.3857142857142857
$ws.Range("D5").Value = 0.3777777777777778
$ws.Range("B6").Value = 0.3958333333333333
$ws.Range("D6").Value = 0.4037037037037037
$ws.Range("B7").Value = 0.5
$ws.Range("C7").Value = 0.7
$ws.Range("D7").Value = 0.5833333333333334
$ws.Range("B8").Value = 0.7
$ws.Range("C8").Value = 0.5
$ws.Range("D8").Value = 0.5833333333333334
$ws.Range("B9").Value = 0.5833333333333334
$ws.Range("C9").Value = 0.5833333333333334
$ws.Range("D9").Value = 0.5833333333333334
$ws.Range("E9").Value = 0.5833333333333334
$ws.Range("B10").Value = 0.6
$ws.Range("C10").Value = 0.6
$ws.Range("D10").Value = 0.5833333333333334
$ws.Range("B11").Value = 0.6166666666666666
$ws.Range("C11").Value = 0.5833333333333334
$ws.Range("D11").Value = 0.5833333333333334
$ws.Range("C12").Value = 0.2
$ws.Range("D12").Value = 0.2857142857142858
$ws.Range("B13").Value = 0.6
$ws.Range("C13").Value = 0.8571428571428571
$ws.Range("D13").Value = 0.7058823529411764
$ws.Range("B15").Value = 0.55
$ws.Range("C15").Value = 0.5285714285714286
$ws.Range("D15").Value = 0.4957983193277311
$ws.Range("B16").Value = 0.5583333333333333
$ws.Range("D16").Value = 0.530812324929972
$ws.Range("B17").Value = 0.2222222222222222
$ws.Range("C17").Value = 0.2
$ws.Range("D17").Value = 0.2105263157894737
$ws.Range("B18").Value = 0.4666666666666667
$ws.Range("C18").Value = 0.5
$ws.Range("D18").Value = 0.4827586206896552
$ws.Range("B19").Value = 0.375
$ws.Range("C19").Value = 0.375
$ws.Range("D19").Value = 0.375
$ws.Range("E19").Value = 0.375
$ws.Range("B20").Value = 0.3444444444444444
$ws.Range("C20").Value = 0.35
$ws.Range("D20").Value = 0.3466424682395645
$ws.Range("B21").Value = 0.3648148148148149
$ws.Range("C21").Value = 0.375
$ws.Range("D21").Value = 0.3693284936479129
$ws.Range("B22").Value = 0.4545454545454545
$ws.Range("C22").Value = 0.5
$ws.Range("D22").Value = 0.4761904761904762
$ws.Range("B23").Value = 0.6153846153846154
$ws.Range("D23").Value = 0.5925925925925927
$ws.Range("B24").Value = 0.5416666666666666
$ws.Range("C24").Value = 0.5416666666666666
$ws.Range("D24").Value = 0.5416666666666666
$ws.Range("E24").Value = 0.5416666666666666
$ws.Range("B25").Value = 0.534965034965035
$ws.Range("C25").Value = 0.5357142857142857
$ws.Range("D25").Value = 0.5343915343915344
$ws.Range("B26").Value = 0.5483682983682984
$ws.Range("C26").Value = 0.5416666666666666
$ws.Range("D26").Value = 0.5440917107583775
